$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple price/volume updates ---
$ws.Cells.Item(2, 4).Value = "55.843.83"
$ws.Cells.Item(2, 5).Value = "  +5.27%  "
$ws.Cells.Item(3, 4).Value = "2.517.95"
$ws.Cells.Item(3, 5).Value = "  +8.34%  "
$ws.Cells.Item(4, 5).Value = "  -0.10%  "
$ws.Cells.Item(5, 4).Value = "485.79"
$ws.Cells.Item(5, 5).Value = "  +11.07%  "
$ws.Cells.Item(6, 4).Value = "142.68"
$ws.Cells.Item(6, 5).Value = "  +17.65%  "
$ws.Cells.Item(7, 5).Value = "  +0.24%  "
$ws.Cells.Item(8, 4).Value = "0.514"
$ws.Cells.Item(8, 5).Value = "  +8.44%  "
$ws.Cells.Item(9, 4).Value = "2.513.53"
$ws.Cells.Item(9, 5).Value = "  +7.61%  "
$ws.Cells.Item(10, 4).Value = "0.0995"
$ws.Cells.Item(10, 5).Value = "  +9.21%  "
$ws.Cells.Item(11, 4).Value = "5.51"
$ws.Cells.Item(11, 5).Value = "  +4.50%  "
$ws.Cells.Item(12, 4).Value = "0.331"
$ws.Cells.Item(12, 5).Value = "  +7.36%  "
$ws.Cells.Item(13, 5).Value = "  +0.68%  "
$ws.Cells.Item(14, 4).Value = "2.947.64"
$ws.Cells.Item(14, 5).Value = "  +8.08%  "
$ws.Cells.Item(15, 4).Value = "55.837.44"
$ws.Cells.Item(15, 5).Value = "  +5.26%  "
$ws.Cells.Item(16, 4).Value = "20.72"
$ws.Cells.Item(16, 5).Value = "  +8.59%  "
$ws.Cells.Item(17, 4).Value = "0.0000139"
$ws.Cells.Item(17, 5).Value = "  +16.08%  "
$ws.Cells.Item(18, 4).Value = "2.515.78"
$ws.Cells.Item(18, 5).Value = "  +7.60%  "
$ws.Cells.Item(19, 4).Value = "4.40"
$ws.Cells.Item(19, 5).Value = "  +10.96%  "
$ws.Cells.Item(20, 4).Value = "321.76"
$ws.Cells.Item(20, 5).Value = "  +6.82%  "
$ws.Cells.Item(21, 4).Value = "10.09"
$ws.Cells.Item(21, 5).Value = "  +10.50%  "
$ws.Cells.Item(22, 5).Value = "  +0.02%  "
$ws.Cells.Item(23, 4).Value = "5.74"
$ws.Cells.Item(23, 5).Value = "  +6.85%  "
$ws.Cells.Item(24, 4).Value = "58.19"
$ws.Cells.Item(24, 5).Value = "  +5.33%  "
$ws.Cells.Item(25, 4).Value = "0.168"
$ws.Cells.Item(25, 5).Value = "  +9.69%  "
$ws.Cells.Item(26, 5).Value = "  +11.29%  "
$ws.Cells.Item(27, 5).Value = "  +0.23%  "
$ws.Cells.Item(28, 4).Value = "2.612.09"
$ws.Cells.Item(28, 5).Value = "  +7.46%  "
$ws.Cells.Item(29, 4).Value = "7.48"
$ws.Cells.Item(29, 5).Value = "  +7.69%  "
$ws.Cells.Item(30, 4).Value = "0.0₃0810"
$ws.Cells.Item(30, 5).Value = "  +15.83%  "
$ws.Cells.Item(31, 4).Value = "1.00"
$ws.Cells.Item(31, 5).Value = "  +0.37%  "
$ws.Cells.Item(32, 4).Value = "149.74"
$ws.Cells.Item(32, 5).Value = "  +4.26%  "
$ws.Cells.Item(33, 4).Value = "18.25"
$ws.Cells.Item(33, 5).Value = "  +6.13%  "
$ws.Cells.Item(34, 5).Value = "  +10.79%  "
$ws.Cells.Item(35, 4).Value = "5.22"
$ws.Cells.Item(35, 5).Value = "  +10.03%  "
$ws.Cells.Item(38, 5).Value = "  +12.96%  "
$ws.Cells.Item(39, 4).Value = "34.29"
$ws.Cells.Item(39, 5).Value = "  +2.74%  "
$ws.Cells.Item(42, 4).Value = "0.0556"
$ws.Cells.Item(42, 5).Value = "  +10.90%  "
$ws.Cells.Item(43, 4).Value = "3.44"
$ws.Cells.Item(43, 5).Value = "  +8.47%  "
$ws.Cells.Item(44, 4).Value = "1.34"
$ws.Cells.Item(44, 5).Value = "  +11.38%  "
$ws.Cells.Item(45, 4).Value = "2.000.12"
$ws.Cells.Item(45, 5).Value = "  +5.33%  "
$ws.Cells.Item(48, 5).Value = "  +9.18%  "
$ws.Cells.Item(49, 4).Value = "254.93"
$ws.Cells.Item(49, 5).Value = "  +35.20%  "
$ws.Cells.Item(50, 5).Value = "  +7.88%  "
$ws.Cells.Item(51, 4).Value = "17.70"
$ws.Cells.Item(51, 5).Value = "  +12.96%  "

# --- Row reordering (coin rank swaps) ---
$ws.Cells.Item(36, 2).Value = "Fetch.AI"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Cells.Item(36, 4).Value = "0.883"
$ws.Cells.Item(36, 5).Value = "  +6.37%  "
$ws.Cells.Item(37, 2).Value = "NEARProtocol"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Cells.Item(37, 4).Value = "3.74"
$ws.Cells.Item(37, 5).Value = "  +5.76%  "
$ws.Cells.Item(40, 2).Value = "Mantle"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Cells.Item(40, 4).Value = "0.617"
$ws.Cells.Item(40, 5).Value = "  +17.45%  "
$ws.Cells.Item(41, 2).Value = "FirstDigitalUSD"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Cells.Item(41, 4).Value = "0.998"
$ws.Cells.Item(41, 5).Value = "  +0.09%  "
$ws.Cells.Item(46, 2).Value = "RenderToken"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(46, 4).Value = "4.70"
$ws.Cells.Item(46, 5).Value = "  +19.16%  "
$ws.Cells.Item(47, 2).Value = "WhiteBITCoin"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Cells.Item(47, 4).Value = "10.17"
$ws.Cells.Item(47, 5).Value = "  -1.20%  "
